$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.891.43'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.60%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.787.91'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.44'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.23%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3925'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.96%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3475'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.31'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.04%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.174'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.97%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07583'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.51%  '

# Row 12
$ws.Range("E12").Value = '  +0.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.75'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.51%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.492'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.68%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.792.44'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.48%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.185'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.57%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001108'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.09%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06742'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.52%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.45'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.37%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.51%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.84'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.01%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.598'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.59%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.892.77'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.60%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.48'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.401'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.57%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.525'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.09%  '

# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.540'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.34%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.36'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.22%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '155.89'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.09%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.83'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.06%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.992.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.29%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.367'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.85%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.009'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.26%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08925'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.93%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '13.19'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.51%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02484'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +5.93%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.524'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.21%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.7013'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.16%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06522'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +2.62%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2248'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.08%  '

# Row 41
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.582'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.94%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.274'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.80%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.607'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.22%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.68'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.31%  '

# Row 45
$ws.Range("E45").Value = '  +0.45%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6435'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.63%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.871'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.71%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.166'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.13%  '

# Row 49
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.69'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.99%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07382'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.66%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.266'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +5.00%  '

